$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Mandatory Y/N" column (was column B) - no fields are mandatory anymore
$ws.Columns.Item(2).Delete()

# Insert a new row at the top for the explanatory banner text
$ws.Rows.Item(1).Insert()

# Style + populate the banner row (bold, white text) before merging so the
# merge carries the formatting across the whole merged range
$ws.Range("A1").Value = "No fields are mandatory, if you don't provide it will search whole project root for .spec or .test by default. But after whole search if any .only is found it will switch its default setting that was searching for .spec or .test and running all to just searching for .only specified test and running."
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.ColorIndex = 2
$ws.Rows.Item(1).RowHeight = 49

$ws.Range("A1:C1").Merge()

# Reset the formatting picked up by the merge on the now-blank B1/C1 cells
$ws.Range("B1:C1").Style = "Normal"

# New rows describing the testDir / testMatch properties
$ws.Range("A3").Value = "testDir"
$ws.Range("B3").Value = "Root"
$ws.Range("C3").Value = "searches for by default in specified and directories inside that recursively"

$ws.Range("A4").Value = "testMatch"
$ws.Range("B4").Value = "Root"
$ws.Range("C4").Value = "specify regex or string pattern to search"
